$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.268.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.56%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.440.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.41%  "

$ws.Range("E4").Value = "  +0.40%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9136"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "274.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.78%  "

$ws.Range("E7").Value = "  -1.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3076"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.27"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.023"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06511"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9985"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.354"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.86%  "

$ws.Range("E14").Value = "  +1.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.044"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.65%  "

$ws.Range("E16").Value = "  -0.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.438.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9292"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05626"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.412"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.237"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.278.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.096"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.591.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.77%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "110.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.926"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8047"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.858"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07667"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.460"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  +0.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.679"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.133"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.18%  "

$ws.Range("E39").Value = "  -1.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01983"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1853"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9264"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.130"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -15.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5217"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.484"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("E46").Value = "  -3.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "117.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5103"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.733"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06409"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.86%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9754"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.55%  "
